# Fix duplicate data and correct reconciliation results.
# - Correct a few TradeType/ProductType values in rows 2,3,5,6
# - Renumber the duplicate TradeID block (rows 17-26, originally T006-T015)
#   to T016-T025 so the sheet holds 25 unique trades (T001-T025).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections on the first block of trades (rows 2-6) ---
$ws.Range("C2").Value = "Static"

$ws.Range("B3").Value = "Swap"
$ws.Range("C3").Value = "Amended"

$ws.Range("C5").Value = "Static"

$ws.Range("C6").Value = "Amended"

# --- Renumber duplicated TradeIDs (rows 17-26): T006-T015 -> T016-T025 ---
# Rows 17-26 mirror rows 7-16 (TradeIDs T006-T015); bump each by +10 so the
# sheet ends up with 25 unique trades T001-T025.
$newIds = @("T016", "T017", "T018", "T019", "T020", "T021", "T022", "T023", "T024", "T025")
for ($i = 0; $i -lt $newIds.Length; $i++) {
    $row = 17 + $i
    $ws.Cells.Item($row, 1).Value = $newIds[$i]
}
